# Update "want to go" counts (column F) that changed between data refreshes.
# Sheet 1 = 展览, Sheet 4 = 全部类型 (per xl/workbook.xml sheet order).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet4 = $wb.Worksheets.Item(4)

# Sheet 1 ("展览") updates
$sheet1.Range("F2").Value = 11421
$sheet1.Range("F3").Value = 10787
$sheet1.Range("F4").Value = 601
$sheet1.Range("F6").Value = 981
$sheet1.Range("F8").Value = 50
$sheet1.Range("F11").Value = 10575
$sheet1.Range("F12").Value = 4080
$sheet1.Range("F14").Value = 2451
$sheet1.Range("F17").Value = 103
$sheet1.Range("F19").Value = 11092
$sheet1.Range("F20").Value = 10850

# Sheet 4 ("全部类型") updates
$sheet4.Range("F2").Value = 11421
$sheet4.Range("F3").Value = 10787
$sheet4.Range("F4").Value = 601
$sheet4.Range("F6").Value = 981
$sheet4.Range("F8").Value = 50
$sheet4.Range("F11").Value = 10575
$sheet4.Range("F12").Value = 4080
$sheet4.Range("F14").Value = 2451
$sheet4.Range("F17").Value = 103
$sheet4.Range("F19").Value = 11092
$sheet4.Range("F20").Value = 10850
